$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to end with a blank, date-formatted row at row 31 (A1:C31).
# The edited workbook keeps that same trailing blank row, just 13 rows further
# down (A1:C44), because 13 new assessments (12 quizzes + 2 extra module-problem
# rows) were inserted above it. Insert blank rows first so that quirky row
# keeps its original formatting (B has the date style, C does not).
for ($i = 0; $i -lt 13; $i++) {
    $ws.Rows.Item(31).Insert()
}

# Pre-create the new shared strings in the same order Excel originally created
# them in, using disposable cells far below the table (column A has no default
# style, so the scratch cells leave no trace once cleared below).
$ws.Cells.Item(200, 1).Value = "Quiz 01"
$ws.Cells.Item(201, 1).Value = "Quiz 02"
$ws.Cells.Item(202, 1).Value = "Quiz 03"
$ws.Cells.Item(203, 1).Value = "Quiz 04"
$ws.Cells.Item(204, 1).Value = "Quiz 05"
$ws.Cells.Item(205, 1).Value = "Quiz 06"
$ws.Cells.Item(206, 1).Value = "Quiz 07"
$ws.Cells.Item(207, 1).Value = "Quiz 08"
$ws.Cells.Item(208, 1).Value = "Quiz 10"
$ws.Cells.Item(209, 1).Value = "Quiz 09"
$ws.Cells.Item(210, 1).Value = "Quiz 11"
$ws.Cells.Item(211, 1).Value = "Quiz 12"
$ws.Cells.Item(212, 1).Value = "Module 10 Programming Problems"
$ws.Cells.Item(213, 1).Value = "Module 12 Programming Problems"

# Write the corrected/expanded assessment schedule (row 1 is the header and is
# unchanged; data now runs through row 43, with the trailing blank row at 44).
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Quiz 01"
$ws.Cells.Item(2, 3).Value = 45301
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Survey"
$ws.Cells.Item(3, 3).Value = 45303
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Quiz 02"
$ws.Cells.Item(4, 3).Value = 45308
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Module 1 Programming Problems"
$ws.Cells.Item(5, 3).Value = 45308
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Quiz 03"
$ws.Cells.Item(6, 3).Value = 45315
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Module 2 Programming Problems"
$ws.Cells.Item(7, 3).Value = 45315
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Programming Project 1"
$ws.Cells.Item(8, 3).Value = 45316
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Quiz 04"
$ws.Cells.Item(9, 3).Value = 45322
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Module 3 Programming Problems"
$ws.Cells.Item(10, 3).Value = 45322
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Programming Project 2"
$ws.Cells.Item(11, 3).Value = 45323
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Quiz 05"
$ws.Cells.Item(12, 3).Value = 45329
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "Module 4 Programming Problems"
$ws.Cells.Item(13, 3).Value = 45329
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "Programming Project 3"
$ws.Cells.Item(14, 3).Value = 45330
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "Midterm 1"
$ws.Cells.Item(15, 3).Value = 45336
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "Module 5 Programming Problems"
$ws.Cells.Item(16, 3).Value = 45337
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "Programming Project 4"
$ws.Cells.Item(17, 3).Value = 45338
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "Quiz 06"
$ws.Cells.Item(18, 3).Value = 45343
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "Module 6 Programming Problems"
$ws.Cells.Item(19, 3).Value = 45343
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "Programming Project 5"
$ws.Cells.Item(20, 3).Value = 45344
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "Quiz 07"
$ws.Cells.Item(21, 3).Value = 45364
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "Module 7 Programming Problems"
$ws.Cells.Item(22, 3).Value = 45364
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "Programming Project 6"
$ws.Cells.Item(23, 3).Value = 45365
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "Quiz 08"
$ws.Cells.Item(24, 3).Value = 45371
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "Module 8 Programming Problems"
$ws.Cells.Item(25, 3).Value = 45371
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "Programming Project 7"
$ws.Cells.Item(26, 3).Value = 45372
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "Midterm 2"
$ws.Cells.Item(27, 3).Value = 45378
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "Module 9 Programming Problems"
$ws.Cells.Item(28, 3).Value = 45379
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = "Quiz 09"
$ws.Cells.Item(29, 3).Value = 45385
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = "Module 10 Programming Problems"
$ws.Cells.Item(30, 3).Value = 45385
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = "Programming Project 8"
$ws.Cells.Item(31, 3).Value = 45386
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = "Quiz 10"
$ws.Cells.Item(32, 3).Value = 45392
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = "Module 11 Programming Problems"
$ws.Cells.Item(33, 3).Value = 45392
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = "Programming Project 9"
$ws.Cells.Item(34, 3).Value = 45393
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = "Quiz 11"
$ws.Cells.Item(35, 3).Value = 45399
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "Module 12 Programming Problems"
$ws.Cells.Item(36, 3).Value = 45399
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "Programming Project 10"
$ws.Cells.Item(37, 3).Value = 45400
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "Midterm 3"
$ws.Cells.Item(38, 3).Value = 45406
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "Module 13 Programming Problems"
$ws.Cells.Item(39, 3).Value = 45407
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "Module 14 Programming Problems"
$ws.Cells.Item(40, 3).Value = 45413
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "Quiz 12"
$ws.Cells.Item(41, 3).Value = 45413
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = "Programming Project 11"
$ws.Cells.Item(42, 3).Value = 45413
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "Final Exam"
$ws.Cells.Item(43, 3).Value = 45415

# Drop the scratch values now that the real cells above reference those strings.
$ws.Cells.Item(200, 1).Value = ""
$ws.Cells.Item(201, 1).Value = ""
$ws.Cells.Item(202, 1).Value = ""
$ws.Cells.Item(203, 1).Value = ""
$ws.Cells.Item(204, 1).Value = ""
$ws.Cells.Item(205, 1).Value = ""
$ws.Cells.Item(206, 1).Value = ""
$ws.Cells.Item(207, 1).Value = ""
$ws.Cells.Item(208, 1).Value = ""
$ws.Cells.Item(209, 1).Value = ""
$ws.Cells.Item(210, 1).Value = ""
$ws.Cells.Item(211, 1).Value = ""
$ws.Cells.Item(212, 1).Value = ""
$ws.Cells.Item(213, 1).Value = ""

# Match the workbooks last-saved selection (the new last row of the table).
$null = $ws.Range("A44").Select()
